# Update cryptos list with latest price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" '57.316.69'
Set-TextCell "E2" '  +0.73%  '

Set-TextCell "D3" '3.056.07'
Set-TextCell "E3" '  +2.67%  '

Set-TextCell "D4" '1.00'
Set-TextCell "E4" '  -0.21%  '

Set-TextCell "D5" '516.48'
Set-TextCell "E5" '  +3.63%  '

Set-TextCell "D6" '141.28'
Set-TextCell "E6" '  +2.84%  '

Set-TextCell "D7" '0.999'
Set-TextCell "E7" '  -0.15%  '

Set-TextCell "E8" '  +2.47%  '

Set-TextCell "D9" '7.24'
Set-TextCell "E9" '  -3.14%  '

Set-TextCell "D10" '0.110'
Set-TextCell "E10" '  +1.95%  '

Set-TextCell "E11" '  +5.49%  '

Set-TextCell "D12" '3.580.09'
Set-TextCell "E12" '  +2.01%  '

Set-TextCell "D13" '0.126'
Set-TextCell "E13" '  -2.12%  '

Set-TextCell "D14" '26.93'
Set-TextCell "E14" '  +4.30%  '

Set-TextCell "D15" '0.0000167'
Set-TextCell "E15" '  +5.54%  '

Set-TextCell "D16" '57.271.33'
Set-TextCell "E16" '  +0.48%  '

Set-TextCell "E17" '  +1.38%  '

Set-TextCell "D18" '3.055.31'
Set-TextCell "E18" '  +2.12%  '

Set-TextCell "D19" '13.44'
Set-TextCell "E19" '  +6.38%  '

Set-TextCell "D20" '8.11'
Set-TextCell "E20" '  +4.25%  '

Set-TextCell "D21" '336.12'
Set-TextCell "E21" '  +4.55%  '

Set-TextCell "E22" '  +0.51%  '

Set-TextCell "D23" '0.509'
Set-TextCell "E23" '  +4.25%  '

Set-TextCell "D24" '65.71'
Set-TextCell "E24" '  +2.88%  '

Set-TextCell "D25" '3.179.66'
Set-TextCell "E25" '  +1.78%  '

Set-TextCell "E26" '  +0.38%  '

Set-TextCell "D27" '0.164'
Set-TextCell "E27" '  +0.97%  '

Set-TextCell "D28" '0.0₃0906'
Set-TextCell "E28" '  +1.46%  '

Set-TextCell "D29" '6.74'
Set-TextCell "E29" '  +2.41%  '

Set-TextCell "D30" '7.21'
Set-TextCell "E30" '  +1.34%  '

Set-TextCell "D31" '1.82'
Set-TextCell "E31" '  +2.31%  '

Set-TextCell "E32" '  +3.27%  '

Set-TextCell "D33" '20.79'
Set-TextCell "E33" '  +3.14%  '

Set-TextCell "D34" '4.74'
Set-TextCell "E34" '  +2.28%  '

Set-TextCell "D35" '153.32'
Set-TextCell "E35" '  -0.15%  '

Set-TextCell "D36" '5.96'
Set-TextCell "E36" '  +3.19%  '

Set-TextCell "D37" '1.28'
Set-TextCell "E37" '  +2.74%  '

Set-TextCell "D38" '25.32'
Set-TextCell "E38" '  +4.91%  '

Set-TextCell "D39" '0.0675'
Set-TextCell "E39" '  +1.47%  '

Set-TextCell "D40" '3.091.06'
Set-TextCell "E40" '  +2.54%  '

Set-TextCell "D41" '37.11'
Set-TextCell "E41" '  -1.18%  '

Set-TextCell "D42" '3.89'
Set-TextCell "E42" '  +3.69%  '

Set-TextCell "D43" '1.00'
Set-TextCell "E43" '  -0.14%  '

Set-TextCell "D44" '0.666'
Set-TextCell "E44" '  +3.65%  '

Set-TextCell "B45" 'Maker'
Set-TextCell "C45" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell "D45" '2.217.89'
Set-TextCell "E45" '  +0.50%  '

Set-TextCell "B46" 'Stacks'
Set-TextCell "C46" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D46" '1.40'
Set-TextCell "E46" '  +1.62%  '

Set-TextCell "D47" '0.968'
Set-TextCell "E47" '  +1.97%  '

Set-TextCell "D48" '6.04'
Set-TextCell "E48" '  +1.95%  '

Set-TextCell "D49" '20.40'
Set-TextCell "E49" '  +6.76%  '

Set-TextCell "D50" '0.0243'
Set-TextCell "E50" '  +3.99%  '

Set-TextCell "E51" '  +12.20%  '

